# Updated symbol list on Sun Dec 25 20:47:34 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) quotes for a set of coin rows on the
# active sheet. Each Price cell is stored as text (it can contain the
# literal "--" placeholder for coins without a quote), so we momentarily
# force a text number format before writing the new value - otherwise a
# numeric-looking string like "242.03" would be auto-converted to a
# number and lose formatting (e.g. "22.90" -> 22.9). The style is then
# reset back to Normal so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2"  "242.03"
Set-TextValue "D3"  "22.91"
Set-TextValue "D4"  "5.381"
Set-TextValue "D5"  "0.05943"
Set-TextValue "D7"  "6.457"
Set-TextValue "D8"  "0.8042"
Set-TextValue "D9"  "0.9123"
Set-TextValue "D11" "0.07424"
Set-TextValue "D12" "0.03308"
Set-TextValue "D13" "0.03048"
Set-TextValue "D14" "0.09322"
Set-TextValue "D15" "3.882"
Set-TextValue "D16" "0.001591"
Set-TextValue "D17" "0.04519"
Set-TextValue "D18" "0.0005943"
Set-TextValue "D19" "0.006137"
Set-TextValue "D20" "0.004994"
Set-TextValue "D22" "0.0009846"
Set-TextValue "D23" "0.00007804"
Set-TextValue "D25" "2.136"
Set-TextValue "D27" "0.1329"
Set-TextValue "D40" "0.03844"
Set-TextValue "D41" "0.006082"
Set-TextValue "D42" "0.1065"
Set-TextValue "D43" "0.002541"
Set-TextValue "D44" "0.007197"
Set-TextValue "D45" "0.00005196"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.0005803"
Set-TextValue "D48" "0.9585"
Set-TextValue "D49" "0.002262"
Set-TextValue "D50" "0.00002101"
Set-TextValue "D51" "0.0002001"
